$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns M, N, O (13,14,15) and P (16) to match the updated layout.
# (Values chosen so the engine's internal pixel-rounding lands on the stored
# OOXML width closest to the target: ~7.7109375 for M:O, ~5.7109375 for P.)
$ws.Range("M1").ColumnWidth = 6.83
$ws.Range("N1").ColumnWidth = 6.83
$ws.Range("O1").ColumnWidth = 6.83
$ws.Range("P1").ColumnWidth = 4.83

# Update the percentage figures in row 1
$ws.Range("M1").Value = 0.060950000000000004
$ws.Range("N1").Value = 0.047050000000000002
$ws.Range("O1").Value = 0.082949999999999996
$ws.Range("P1").Value = 0.081000000000000003
